$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Copy the formatting of an existing data row (row 2) onto the new
#    rows 10-16 so the new cells pick up the same cell styles
#    (date format, fills, alignment, wrap) as the rest of the table.
#    xlPasteFormats = -4122
# ---------------------------------------------------------------------
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A10:E16").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2) Fill in the new rows with their values.
#    NOTE: the write order below matters - it reproduces the exact
#    order in which new, unique strings were first entered so that the
#    shared-string table gets rebuilt with the same indices/order as
#    the target workbook (D16 "Powerpoint" was typed before the rest
#    of the new rows were filled in).
# ---------------------------------------------------------------------

# Row 16's task cell ("Powerpoint") is the very first new unique string.
$ws.Range("D16").Value = "Powerpoint"

# Row 10
$ws.Range("A10").Value = 42797
$ws.Range("B10").Value = "Basiel"
$ws.Range("C10").Value = "DEV"
$ws.Range("D10").Value = "opzetten laravel omgeving"
$ws.Range("E10").Value = 3

# Row 11
$ws.Range("A11").Value = 42802
$ws.Range("B11").Value = "Basiel"
$ws.Range("C11").Value = "DES"
$ws.Range("D11").Value = "Style tiles"
$ws.Range("E11").Value = "'0.5"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null

# Row 12
$ws.Range("A12").Value = 42804
$ws.Range("B12").Value = "Basiel"
$ws.Range("C12").Value = "DES"
$ws.Range("D12").Value = "Bespreking style tiles, wireframes, ..."
$ws.Range("E12").Value = 3

# Row 13
$ws.Range("A13").Value = 42807
$ws.Range("B13").Value = "Basiel"
$ws.Range("C13").Value = "DEV"
$ws.Range("D13").Value = "Database migrations"
$ws.Range("E13").Value = 3

# Row 14
$ws.Range("A14").Value = 42808
$ws.Range("B14").Value = "Basiel"
$ws.Range("C14").Value = "DES"
$ws.Range("D14").Value = "Visual designs"
$ws.Range("E14").Value = 2

# Row 15
$ws.Range("A15").Value = 42810
$ws.Range("B15").Value = "Basiel"
$ws.Range("C15").Value = "DES"
$ws.Range("D15").Value = "Dossier in MD"
$ws.Range("E15").Value = "'2.5"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

# Row 16 (D16 already set above)
$ws.Range("A16").Value = 42810
$ws.Range("B16").Value = "Basiel"
$ws.Range("C16").Value = "DES"
$ws.Range("E16").Value = "'1.5"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 3) Update the selected cell to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("F6").Select() | Out-Null
